$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells for the season record: Wins / Losses / Ties.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the rest of the header row (bold, boxed,
# centered) by copying the existing AC1 header cell's formats onto the
# new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins=70, Losses=92, Ties=0) for every
# player row, 2 through 56.
$ws.Range("AD2:AD56").Value = 70
$ws.Range("AE2:AE56").Value = 92
$ws.Range("AF2:AF56").Value = 0
